$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 605, shifting existing rows 605-649 down to 606-650.
$ws.Rows("605:605").Insert()

# Populate the newly inserted row 605 with the new record.
$ws.Range("A605").Value = 6
$ws.Range("B605").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C605").Value = "Metropolitana"
$ws.Range("D605").Value = 45041
$ws.Range("E605").Value = 13
$ws.Range("F605").Value = "Fruta"
$ws.Range("G605").Value = 100101
$ws.Range("H605").Value = "Berries"
$ws.Range("I605").Value = 100101001
$ws.Range("J605").Value = "Arándano (blue)"
$ws.Range("K605").Value = "Sin especificar"
$ws.Range("L605").Value = "Primera"
$ws.Range("M605").Value = 100
$ws.Range("N605").Value = 10000
$ws.Range("O605").Value = 10000
$ws.Range("P605").Value = 10000
$ws.Range("Q605").Value = "$/bandeja 2 kilos"
$ws.Range("R605").Value = "Provincia de Curicó"
$ws.Range("S605").Value = 5000
$ws.Range("T605").Value = 2
